$p = $ppt.ActivePresentation

foreach ($idx in 14,15,16) {
    $slide = $p.Slides.Item($idx)
    $tbl = $slide.Shapes.Item(1).Table
    $tbl.ApplyStyle("{5DA2F8D5-C1A0-4302-B0CC-28C5FCE34058}")
}
